# Update countries & provincias Spain
# Applies the daily COVID-19 data refresh to the "Pais" worksheet:
#  - numeric updates for several countries (Rusia, India, Barein, Armenia,
#    Moldavia, Chequia, Estonia, Taiwan)
#  - alphabetical re-sort swap of Fiyi/Dominica and of
#    Islas Turcas y Caicos/Santa Sede (including their row data)
#  - refresh of the "last updated" timestamp banner in A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rusia (row 6): Casos totales, Nuevos casos, Casos activos, Recuperados,
# Muertes hoy, Muertes
$ws.Range("B6").Value = 561091
$ws.Range("C6").Value = 7790
$ws.Range("D6").Value = 313963
$ws.Range("E6").Value = 239468
$ws.Range("G6").Value = 182
$ws.Range("H6").Value = 7660

# India (row 7): Casos activos, Recuperados
$ws.Range("D7").Value = 194565
$ws.Range("E7").Value = 161126

# Barein (row 50): Recuperados, Muertes hoy, Muertes
$ws.Range("E50").Value = 5725
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 51

# Armenia (row 52): Casos totales, Nuevos casos, Casos activos,
# Recuperados, Muertes hoy, Muertes
$ws.Range("B52").Value = 18698
$ws.Range("C52").Value = 665
$ws.Range("D52").Value = 7560
$ws.Range("E52").Value = 10829
$ws.Range("G52").Value = 7
$ws.Range("H52").Value = 309

# Moldavia (row 57): Casos activos, Recuperados, Muertes hoy, Muertes
$ws.Range("D57").Value = 7252
$ws.Range("E57").Value = 5042
$ws.Range("G57").Value = 5
$ws.Range("H57").Value = 438

# Chequia (row 66): Casos totales, Nuevos casos, Recuperados
$ws.Range("B66").Value = 10176
$ws.Range("C66").Value = 14
$ws.Range("E66").Value = 2444

# Estonia (row 103): Casos activos, Recuperados
$ws.Range("D103").Value = 1748
$ws.Range("E103").Value = 160

# Taiwan (row 154): Casos totales, Nuevos casos, Recuperados
$ws.Range("B154").Value = 446
$ws.Range("C154").Value = 1
$ws.Range("E154").Value = 5

# Re-sort: Fiyi/Dominica swap places (row 202/203 numbers are identical,
# only the country labels swap)
$ws.Range("A202").Value = "Dominica"
$ws.Range("A203").Value = "Fiyi"

# Re-sort: Islas Turcas y Caicos/Santa Sede swap places (labels and data
# both move between row 208 and row 209)
$ws.Range("A208").Value = "Santa Sede"
$ws.Range("D208").Value = 12
$ws.Range("H208").Value = 0

$ws.Range("A209").Value = "Islas Turcas y Caicos"
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 1

# Refresh "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 18 de Junio de 2020 a las 10:13"
